# Apply updated crypto price/volume data (GitHub Actions refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    # Force the cell to stay text (matches source inlineStr cells) even
    # when the string looks numeric (e.g. "0.9999", "327.15"), then drop
    # back to the default "Normal" style so no stray number format sticks.
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "29.536.10"
$ws.Range("E2").Value = "  +1.33%  "
Set-TextValue $ws.Range("D3") "1.918.04"
$ws.Range("E3").Value = "  +0.26%  "
Set-TextValue $ws.Range("D4") "0.9999"
$ws.Range("E4").Value = "  -0.01%  "
Set-TextValue $ws.Range("D5") "327.15"
$ws.Range("E5").Value = "  -2.07%  "
Set-TextValue $ws.Range("D6") "0.9999"
$ws.Range("E6").Value = "  +0.08%  "
Set-TextValue $ws.Range("D7") "0.4803"
$ws.Range("E7").Value = "  +3.25%  "
Set-TextValue $ws.Range("D8") "0.4106"
$ws.Range("E8").Value = "  +0.04%  "
Set-TextValue $ws.Range("D9") "47.68"
$ws.Range("E9").Value = "  -0.02%  "
Set-TextValue $ws.Range("D10") "0.08054"
$ws.Range("E10").Value = "  +0.33%  "
Set-TextValue $ws.Range("D11") "1.013"
$ws.Range("E11").Value = "  +0.10%  "
Set-TextValue $ws.Range("D12") "22.44"
$ws.Range("E12").Value = "  +2.14%  "
Set-TextValue $ws.Range("D13") "1.903.57"
$ws.Range("E13").Value = "  -0.03%  "
Set-TextValue $ws.Range("D14") "5.958"
$ws.Range("E14").Value = "  -0.16%  "
Set-TextValue $ws.Range("D15") "7.163"
$ws.Range("E15").Value = "  +0.62%  "
Set-TextValue $ws.Range("D16") "89.70"
$ws.Range("E16").Value = "  +0.31%  "
Set-TextValue $ws.Range("D17") "0.9998"
$ws.Range("E17").Value = "  -0.03%  "
Set-TextValue $ws.Range("D20") "17.78"
$ws.Range("E20").Value = "  +1.02%  "
Set-TextValue $ws.Range("D21") "1.001"
$ws.Range("E21").Value = "  +0.01%  "
Set-TextValue $ws.Range("D22") "29.545.09"
$ws.Range("E22").Value = "  +1.15%  "
Set-TextValue $ws.Range("D23") "5.556"
$ws.Range("E23").Value = "  +1.82%  "
$ws.Range("E24").Value = "  +1.72%  "
$ws.Range("E25").Value = "  -1.52%  "
Set-TextValue $ws.Range("D26") "2.133.39"
$ws.Range("E26").Value = "  +0.29%  "
Set-TextValue $ws.Range("D27") "153.57"
$ws.Range("E27").Value = "  -2.54%  "
Set-TextValue $ws.Range("D28") "19.89"
$ws.Range("E28").Value = "  +0.54%  "
Set-TextValue $ws.Range("D29") "5.820"
$ws.Range("E29").Value = "  +6.80%  "
Set-TextValue $ws.Range("D30") "2.142"
$ws.Range("E30").Value = "  +0.87%  "
Set-TextValue $ws.Range("D31") "118.04"
$ws.Range("E31").Value = "  -0.52%  "
Set-TextValue $ws.Range("D32") "1.060"
$ws.Range("E32").Value = "  +6.90%  "
Set-TextValue $ws.Range("D33") "0.09538"
$ws.Range("E33").Value = "  +1.03%  "
Set-TextValue $ws.Range("D34") "1.427"
$ws.Range("E34").Value = "  -1.07%  "
$ws.Range("E35").Value = "  -0.39%  "
Set-TextValue $ws.Range("D36") "5.393"
$ws.Range("E36").Value = "  +1.15%  "
Set-TextValue $ws.Range("D37") "0.06107"
$ws.Range("E37").Value = "  -0.11%  "
Set-TextValue $ws.Range("D38") "0.02258"
$ws.Range("E38").Value = "  +0.11%  "
Set-TextValue $ws.Range("D39") "8.353"
$ws.Range("E39").Value = "  -0.66%  "
Set-TextValue $ws.Range("D40") "1.173"
$ws.Range("E40").Value = "  -0.65%  "
Set-TextValue $ws.Range("D41") "0.5894"
$ws.Range("E41").Value = "  +0.99%  "
$ws.Range("E46").Value = "  +3.07%  "
Set-TextValue $ws.Range("D49") "1.935"
$ws.Range("E49").Value = "  +0.39%  "
Set-TextValue $ws.Range("D50") "113.59"
$ws.Range("E50").Value = "  +1.49%  "
Set-TextValue $ws.Range("D51") "44.88"
$ws.Range("E51").Value = "  -5.80%  "

$ws.Range("B18").Value = "ShibaInu"
$ws.Range("C18").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
Set-TextValue $ws.Range("D18") "0.00001033"
$ws.Range("E18").Value = "  -0.26%  "

$ws.Range("B19").Value = "TRON"
$ws.Range("C19").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
Set-TextValue $ws.Range("D19") "0.06605"
$ws.Range("E19").Value = "  +0.31%  "

$ws.Range("B42").Value = "RenderToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue $ws.Range("D42") "2.533"
$ws.Range("E42").Value = "  +6.56%  "

$ws.Range("B43").Value = "Algorand"
$ws.Range("C43").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
Set-TextValue $ws.Range("D43") "0.1849"
$ws.Range("E43").Value = "  +0.89%  "

$ws.Range("B44").Value = "Aptos"
$ws.Range("C44").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
Set-TextValue $ws.Range("D44") "10.16"
$ws.Range("E44").Value = "  -0.56%  "

$ws.Range("B45").Value = "Cronos"
$ws.Range("C45").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextValue $ws.Range("D45") "0.07973"
$ws.Range("E45").Value = "  +12.56%  "

$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue $ws.Range("D47") "12.22"
$ws.Range("E47").Value = "  +0.50%  "

$ws.Range("B48").Value = "Decentraland"
$ws.Range("C48").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
Set-TextValue $ws.Range("D48") "0.5557"
$ws.Range("E48").Value = "  +0.41%  "

